# "Generate Report for Archive"
#
# The localization-status report previously showed the placeholder status
# text "Ready for handoff" for the two sample rows on every status-bearing
# column (the "zh-cn"/"de-de" status columns on the Overview sheet, and the
# "Status" column on each of the per-locale sheets). The archived report
# now reflects that those files are still "In Translation", so every cell
# that used to read "Ready for handoff" is updated to "In Translation".
#
# Updating the text shrinks the longest string in those columns, so the
# columns that only ever held this status text are narrowed to match
# (AutoFit-style resize driven by the new, shorter content).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn (col E) / de-de (col F) status cells ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

# --- zh-cn sheet: Status column (col C) ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"

# --- de-de sheet: Status column (col C) ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"

# Narrow the now-shorter status columns to fit the new text.
$ws1.Columns.Item(5).ColumnWidth = 12.5   # Overview!E (zh-cn status)
$ws1.Columns.Item(6).ColumnWidth = 12.5   # Overview!F (de-de status)
$ws2.Columns.Item(3).ColumnWidth = 12.5   # zh-cn!C (Status)
$ws3.Columns.Item(3).ColumnWidth = 12.5   # de-de!C (Status)
